$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2:R3 down to whole numbers
$ws.Range("Q2").Value = 549104
$ws.Range("R2").Value = 6531856
$ws.Range("Q3").Value = 549104
$ws.Range("R3").Value = 6531856

# Clear the Starttid (Z) and Sluttid (AB) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
